$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Replace the RMA-ZYC7* test-run identifiers with a freshly generated
# RMA-M7GG* run (new Test plan execution values / Salesforce Ids), row by row.
$ws.Range("E2").Value2 = "RMA-M7GG-001"
$ws.Range("F2").Value2 = "RMA-M7GG-1-1"
$ws.Range("J2").Value2 = "a7s5f000000xNYDAA2"

$ws.Range("E3").Value2 = "RMA-M7GG-002"
$ws.Range("F3").Value2 = "RMA-M7GG-1-2"
$ws.Range("J3").Value2 = "a7s5f000000xNYEAA2"

$ws.Range("E4").Value2 = "RMA-M7GG-003"
$ws.Range("F4").Value2 = "RMA-M7GG-1-3"
$ws.Range("J4").Value2 = "a7s5f000000xNYFAA2"

# The new Ids are the same character length as the old ones, but Excel's
# "best fit" column-width recalculation still shifts slightly because of
# glyph-width differences, so nudge the affected columns to match.
$ws.Columns.Item(5).ColumnWidth = 15.0
$ws.Columns.Item(6).ColumnWidth = 14.5
$ws.Columns.Item(10).ColumnWidth = 20.166666666666668
